$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 7 on the Overview sheet corresponds to file 910bdb67-0a00-40aa-bd98-949a431a4405.md
# Status moves from "Ready for handoff" to "In Translation" for both zh-cn and de-de columns.
$wsOverview.Range("E7").Value = "In Translation"
$wsOverview.Range("F7").Value = "In Translation"

# Row 7 on the zh-cn sheet: Status column C
$wsZhCn.Range("C7").Value = "In Translation"

# Row 7 on the de-de sheet: Status column C
$wsDeDe.Range("C7").Value = "In Translation"
